# Update cryptos list values per upstream diff (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '37.294.34'
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '  -1.39%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.041.77'
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '  -2.30%  '
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '  +0.36%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '228.38'
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '  -2.31%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.611'
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '  -2.22%  '
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '  +0.08%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '56.02'
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '  -4.44%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.383'
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '  -3.05%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.0808'
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '  +2.88%  '
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '  -1.96%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '2.345.75'
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '  -2.18%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '14.50'
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '  -4.18%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '20.44'
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '  -4.59%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.749'
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '  -4.17%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '5.25'
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '  -2.32%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '2.043.08'
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '  -2.50%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '37.144.90'
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '  -1.81%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '5.94'
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '  -3.39%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '69.48'
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '  -2.58%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.0$([char]0x2083)0837"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '  -0.22%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '225.11'
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '  -2.26%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '  +0.04%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.35'
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '  -1.94%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.26'
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '  -6.34%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '9.46'
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '  -3.81%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '167.84'
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = '  -2.36%  '
$cell = $ws.Range("B28")
$cell.NumberFormat = "@"
$cell.Value = 'Kaspa'
$cell = $ws.Range("C28")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '0.128'
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = '  -7.00%  '
$cell = $ws.Range("B29")
$cell.NumberFormat = "@"
$cell.Value = 'ImmutableX'
$cell = $ws.Range("C29")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.38'
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = '  -1.27%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '18.85'
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = '  -3.59%  '
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = '  -3.25%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '4.48'
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = '  -5.13%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '4.55'
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = '  -2.85%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.0608'
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = '  -4.11%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '2.38'
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = '  -4.75%  '
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = '  -0.05%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = '  +0.06%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '3.18'
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '  -6.73%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '5.34'
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '  -0.45%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.0219'
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '  -7.94%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '2.89'
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '  -1.25%  '
$cell = $ws.Range("B42")
$cell.NumberFormat = "@"
$cell.Value = 'Maker'
$cell = $ws.Range("C42")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '1.477.32'
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '  +1.55%  '
$cell = $ws.Range("B43")
$cell.NumberFormat = "@"
$cell.Value = 'InjectiveProtocol'
$cell = $ws.Range("C43")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '16.76'
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '  -0.52%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.0936'
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '  -4.00%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '95.12'
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '  -7.28%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '1.14'
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '  -1.47%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.01'
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '  -5.29%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '7.09'
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '  -3.44%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '2.90'
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '  -2.60%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '2.233.22'
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '  -2.10%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '3.60'
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '  -14.88%  '

Write-Output "Applied crypto updates"
